$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "modified" column header
$ws.Range("D1").Value = "modified"

# Mark the "yes" row (row 3) and "y" row (row 5) as modified with an "x"
$ws.Range("D3").Value = "x"
$ws.Range("D5").Value = "x"

# Column C got narrower (resized by the author)
$ws.Columns.Item(3).ColumnWidth = 7.33

# Final selection left on D4
$ws.Range("D4").Select()
